$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Value)
    $ws.Range("ZZ1").Value = "'" + $Value
    $ws.Range("ZZ1").Copy()
    $ws.Range($Cell).PasteSpecial(-4163)
}

Set-TextValue "D2" "311.13"
Set-TextValue "E2" "0.89%"
Set-TextValue "E3" "1.97%"
Set-TextValue "D4" "5.169"
Set-TextValue "E4" "1.18%"
Set-TextValue "D5" "0.08151"
Set-TextValue "E5" "0.26%"
Set-TextValue "D6" "2.005"
Set-TextValue "E6" "2.06%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue "D7" "4.241"
Set-TextValue "E7" "0.77%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
Set-TextValue "D8" "8.146"
Set-TextValue "E8" "2.72%"
Set-TextValue "D9" "0.9265"
Set-TextValue "E9" "-0.25%"
Set-TextValue "D10" "0.1386"
Set-TextValue "E10" "-4.60%"
Set-TextValue "D11" "0.1937"
Set-TextValue "E11" "-1.07%"
Set-TextValue "D12" "0.09066"
Set-TextValue "E12" "-0.83%"
Set-TextValue "D13" "0.03516"
Set-TextValue "E13" "0.21%"
Set-TextValue "D14" "0.09819"
Set-TextValue "E14" "-0.03%"
Set-TextValue "D15" "0.001388"
Set-TextValue "E15" "-1.11%"
Set-TextValue "D16" "0.006100"
Set-TextValue "E16" "0.32%"
Set-TextValue "D17" "3.676"
Set-TextValue "E17" "0.85%"
Set-TextValue "D19" "0.3457"
Set-TextValue "E19" "0.31%"
Set-TextValue "E20" "2.87%"
Set-TextValue "D21" "4.654"
Set-TextValue "E21" "-2.89%"
Set-TextValue "D22" "0.2423"
Set-TextValue "E22" "-1.24%"
Set-TextValue "D23" "0.04376"
Set-TextValue "E23" "-1.50%"
Set-TextValue "D24" "0.001230"
Set-TextValue "E24" "0.96%"
Set-TextValue "D25" "0.004874"
Set-TextValue "E25" "0.68%"
Set-TextValue "D26" "0.0001299"
Set-TextValue "E26" "-0.15%"
Set-TextValue "D27" "0.0003998"
Set-TextValue "E27" "-10.11%"
Set-TextValue "D39" "0.02140"
Set-TextValue "E39" "1.76%"
Set-TextValue "D40" "0.05202"
Set-TextValue "E40" "1.22%"
Set-TextValue "D41" "0.007431"
Set-TextValue "E41" "-0.72%"
Set-TextValue "D42" "0.009834"
Set-TextValue "E42" "-3.03%"
Set-TextValue "D43" "0.1368"
Set-TextValue "E43" "0.39%"
Set-TextValue "D44" "0.002129"
Set-TextValue "E44" "-0.61%"
Set-TextValue "D45" "0.009873"
Set-TextValue "E45" "-3.58%"
Set-TextValue "D46" "0.00006404"
Set-TextValue "E46" "1.88%"
Set-TextValue "D47" "0.00000000750"
Set-TextValue "E47" "-0.15%"
$ws.Range("B48").Value = "BOLO"
$ws.Range("C48").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-TextValue "D48" "0.002753"
Set-TextValue "E48" "-9.99%"
$ws.Range("B49").Value = "CoinbaseStockToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
Set-TextValue "D49" "0.0009995"
Set-TextValue "E49" "-37.59%"
Set-TextValue "D50" "0.00002099"
Set-TextValue "E50" "-0.15%"
Set-TextValue "D51" "0.0001999"
Set-TextValue "E51" "-0.15%"

$ws.Range("ZZ1").Clear()
